{"js": "// Find the \"Literature Review: \" heading paragraph and, right after its\n// existing bold label run, add a new (non-bold) sentence of body text.\n// Also insert a brand-new empty paragraph directly after this paragraph\n// (mirroring the blank \"spacer\" paragraph that already followed it),\n// pushing the original spacer paragraph further down.\nconst searchResults = context.document.body.search(\"Literature Review:\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the 'Literature Review:' paragraph\");\n}\n\nconst headingParagraph = searchResults.items[0].paragraphs.getFirst();\nconst headingRange = headingParagraph.getRange();\n\nfunction wrapAsPackageOoxml(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n    bodyXml +\n    '</w:body></w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst rFonts = '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>';\n\nconst updatedHeadingParagraph =\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:spacing w:line=\"480\" w:lineRule=\"auto\"/>' +\n      '<w:rPr>' + rFonts + '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n      '<w:rPr>' + rFonts + '<w:b/><w:bCs/></w:rPr>' +\n      '<w:t xml:space=\"preserve\">Literature Review: </w:t>' +\n    '</w:r>' +\n    '<w:r>' +\n      '<w:rPr>' + rFonts + '</w:rPr>' +\n      '<w:t xml:space=\"preserve\">Our final dashboard drew inspiration from several ideas we found across literature of previous sports analytics projects. </w:t>' +\n    '</w:r>' +\n  '</w:p>';\n\nconst newBlankParagraph =\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:spacing w:line=\"480\" w:lineRule=\"auto\"/>' +\n      '<w:rPr>' + rFonts + '<w:b/><w:bCs/></w:rPr>' +\n    '</w:pPr>' +\n  '</w:p>';\n\nheadingRange.insertOoxml(\n  wrapAsPackageOoxml(updatedHeadingParagraph + newBlankParagraph),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Find the \"Literature Review: \" heading paragraph and, right after its\n# existing bold label run, add a new (non-bold) sentence of body text.\n# Also insert a brand-new empty paragraph directly after this paragraph\n# (mirroring the blank \"spacer\" paragraph that already followed it),\n# pushing the original spacer paragraph further down.\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"Literature Review:\")\nif (-not $found) {\n    throw \"Could not find the 'Literature Review:' paragraph\"\n}\n\n# Expand the collapsed found-range out to the whole paragraph (wdParagraph = 4)\n# so we can replace its OOXML wholesale with precise formatting.\n[void]$searchRange.Expand(4)\n\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n$rFonts = '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>'\n\n$updatedHeadingParagraph = '<w:p ' + $wNs + '>' +\n    '<w:pPr>' +\n      '<w:spacing w:line=\"480\" w:lineRule=\"auto\"/>' +\n      '<w:rPr>' + $rFonts + '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n      '<w:rPr>' + $rFonts + '<w:b/><w:bCs/></w:rPr>' +\n      '<w:t xml:space=\"preserve\">Literature Review: </w:t>' +\n    '</w:r>' +\n    '<w:r>' +\n      '<w:rPr>' + $rFonts + '</w:rPr>' +\n      '<w:t xml:space=\"preserve\">Our final dashboard drew inspiration from several ideas we found across literature of previous sports analytics projects. </w:t>' +\n    '</w:r>' +\n  '</w:p>'\n\n$newBlankParagraph = '<w:p ' + $wNs + '>' +\n    '<w:pPr>' +\n      '<w:spacing w:line=\"480\" w:lineRule=\"auto\"/>' +\n      '<w:rPr>' + $rFonts + '<w:b/><w:bCs/></w:rPr>' +\n    '</w:pPr>' +\n  '</w:p>'\n\n$searchRange.InsertXML($updatedHeadingParagraph + $newBlankParagraph)\n"}
